$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Extension Payments" Tax Type row (row 3) is excluded from execution
# by flipping its Execute flag from "Y" to "DONOTRUN".
$ws.Range("C3").Value2 = "DONOTRUN"

# Refresh the recorded run timestamps for the rows that still execute.
$ws.Range("B2").Value2 = "Thu Dec 07 21:43:36 EST 2023"
$ws.Range("B4").Value2 = "Thu Dec 07 21:43:49 EST 2023"

# Widen the Execute column so the longer "DONOTRUN" value is fully visible.
$ws.Columns.Item(3).ColumnWidth = 20

# Leave the cell selection resting on the Execute cell that was edited.
$ws.Range("C3").Select()
